# Scheduled market-price refresh: updates currentAveragePrice* / Leve*Price* /
# LeveProfit* columns (H:N) for a number of leve rows across all job sheets,
# reflecting newly-fetched Universalis market data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 94.59999999999999
$ws.Range("I33").Value = 94.59999999999999
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 94.59999999999999
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 134.4
$ws.Range("N33").ClearContents()

$ws.Range("H41").Value = 340.14285
$ws.Range("I41").Value = 196.5
$ws.Range("J41").Value = 531.6667
$ws.Range("K41").Value = 196.5
$ws.Range("L41").Value = 531.6667
$ws.Range("M41").Value = 243.5
$ws.Range("N41").Value = -1411.6667

$ws.Range("H62").Value = 2622.7
$ws.Range("I62").Value = 1959
$ws.Range("K62").Value = 1959
$ws.Range("M62").Value = -1335

$ws.Range("H65").Value = 2622.7
$ws.Range("I65").Value = 1959
$ws.Range("K65").Value = 9795
$ws.Range("M65").Value = -6675

$ws.Range("H137").Value = 48398.816
$ws.Range("I137").Value = 3164.8667
$ws.Range("K137").Value = 9494.6001
$ws.Range("M137").Value = -6944.6001

$ws.Range("H138").Value = 1359.5269
$ws.Range("I138").Value = 549.86664
$ws.Range("J138").Value = 2118.5833
$ws.Range("K138").Value = 1649.59992
$ws.Range("L138").Value = 6355.749899999999
$ws.Range("M138").Value = 3490.40008
$ws.Range("N138").Value = -16635.7499

$ws.Range("H141").Value = 2685
$ws.Range("I141").Value = 2299.1667
$ws.Range("J141").Value = 5000
$ws.Range("K141").Value = 6897.500100000001
$ws.Range("L141").Value = 15000
$ws.Range("M141").Value = -1717.500100000001
$ws.Range("N141").Value = -25360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 14667.205
$ws.Range("I132").Value = 1689.8064
$ws.Range("J132").Value = 64954.625
$ws.Range("K132").Value = 5069.4192
$ws.Range("L132").Value = 194863.875
$ws.Range("M132").Value = -2539.4192
$ws.Range("N132").Value = -199923.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4672.273
$ws.Range("I105").Value = 5681
$ws.Range("J105").Value = 3831.6667
$ws.Range("K105").Value = 5681
$ws.Range("L105").Value = 3831.6667
$ws.Range("M105").Value = -3934
$ws.Range("N105").Value = -7325.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12621.561
$ws.Range("I31").Value = 19651.682
$ws.Range("K31").Value = 19651.682
$ws.Range("M31").Value = -19356.682

$ws.Range("H34").Value = 12621.561
$ws.Range("I34").Value = 19651.682
$ws.Range("K34").Value = 19651.682
$ws.Range("M34").Value = -19449.682

$ws.Range("H86").Value = 6418331
$ws.Range("I86").Value = 2841.1428
$ws.Range("J86").Value = 13903070
$ws.Range("K86").Value = 2841.1428
$ws.Range("L86").Value = 13903070
$ws.Range("M86").Value = -1718.1428
$ws.Range("N86").Value = -13905316

$ws.Range("H89").Value = 6418331
$ws.Range("I89").Value = 2841.1428
$ws.Range("J89").Value = 13903070
$ws.Range("K89").Value = 14205.714
$ws.Range("L89").Value = 69515350
$ws.Range("M89").Value = -8589.714
$ws.Range("N89").Value = -69526582

$ws.Range("H122").Value = 1139.1875
$ws.Range("I122").Value = 1242.8572
$ws.Range("J122").Value = 1058.5555
$ws.Range("K122").Value = 3728.5716
$ws.Range("L122").Value = 3175.6665
$ws.Range("M122").Value = -1278.5716
$ws.Range("N122").Value = -8075.666499999999

$ws.Range("H134").Value = 1074.6774
$ws.Range("I134").Value = 781.5135
$ws.Range("K134").Value = 2344.5405
$ws.Range("M134").Value = 190.4594999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 419.16666
$ws.Range("I26").Value = 171.66667
$ws.Range("K26").Value = 515.00001
$ws.Range("M26").Value = -227.00001

$ws.Range("H107").Value = 7953.231
$ws.Range("I107").Value = 14498
$ws.Range("J107").Value = 317.66666
$ws.Range("K107").Value = 43494
$ws.Range("L107").Value = 952.9999799999999
$ws.Range("M107").Value = -41574
$ws.Range("N107").Value = -4792.99998

$ws.Range("H113").Value = 562.2353000000001
$ws.Range("I113").Value = 491.8
$ws.Range("J113").Value = 591.5833
$ws.Range("K113").Value = 1475.4
$ws.Range("L113").Value = 1774.7499
$ws.Range("M113").Value = 694.5999999999999
$ws.Range("N113").Value = -6114.7499

$ws.Range("H114").Value = 1074.8125
$ws.Range("I114").Value = 1428.2858
$ws.Range("J114").Value = 799.8889
$ws.Range("K114").Value = 4284.857400000001
$ws.Range("L114").Value = 2399.6667
$ws.Range("M114").Value = -1030.857400000001
$ws.Range("N114").Value = -8907.6667

$ws.Range("H131").Value = 795.71
$ws.Range("J131").Value = 796.6767599999999
$ws.Range("L131").Value = 2390.03028
$ws.Range("N131").Value = -12470.03028

$ws.Range("H132").Value = 1719.8
$ws.Range("I132").Value = 1699.6666
$ws.Range("J132").Value = 1750
$ws.Range("K132").Value = 15296.9994
$ws.Range("L132").Value = 15750
$ws.Range("M132").Value = -12766.9994
$ws.Range("N132").Value = -20810

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 638.4211
$ws.Range("I97").Value = 542.2
$ws.Range("K97").Value = 542.2
$ws.Range("M97").Value = -46.20000000000005

$ws.Range("H102").Value = 21741198
$ws.Range("I102").Value = 26318102
$ws.Range("J102").Value = 903.5
$ws.Range("K102").Value = 26318102
$ws.Range("L102").Value = 903.5
$ws.Range("M102").Value = -26316480
$ws.Range("N102").Value = -4147.5

$ws.Range("H122").Value = 53334344
$ws.Range("J122").Value = 90909830
$ws.Range("L122").Value = 272729490
$ws.Range("N122").Value = -272734390

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4029.4285
$ws.Range("I7").Value = 3747.077
$ws.Range("J7").Value = 7700
$ws.Range("K7").Value = 3747.077
$ws.Range("L7").Value = 7700
$ws.Range("M7").Value = -3635.077
$ws.Range("N7").Value = -7924

$ws.Range("H40").Value = 4069.0286
$ws.Range("J40").Value = 5015.524
$ws.Range("L40").Value = 5015.524
$ws.Range("N40").Value = -5287.524

$ws.Range("H61").Value = 4595.75
$ws.Range("I61").Value = 2448.0908
$ws.Range("J61").Value = 9320.6
$ws.Range("K61").Value = 2448.0908
$ws.Range("L61").Value = 9320.6
$ws.Range("M61").Value = -2246.0908
$ws.Range("N61").Value = -9724.6

$ws.Range("H113").Value = 4595.75
$ws.Range("I113").Value = 2448.0908
$ws.Range("J113").Value = 9320.6
$ws.Range("K113").Value = 2448.0908
$ws.Range("L113").Value = 9320.6
$ws.Range("M113").Value = -278.0907999999999
$ws.Range("N113").Value = -13660.6

$ws.Range("H126").Value = 4029.4285
$ws.Range("I126").Value = 3747.077
$ws.Range("J126").Value = 7700
$ws.Range("K126").Value = 11241.231
$ws.Range("L126").Value = 23100
$ws.Range("M126").Value = -8771.231
$ws.Range("N126").Value = -28040

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2271.1428
$ws.Range("I132").Value = 1180
$ws.Range("K132").Value = 3540
$ws.Range("M132").Value = -1010
